# Added search box for events w/ autocomplete
# (underlying data change: add a "Type" column to the "Event Destinations"
#  sheet, marking every row as a "Destination", mirroring the existing
#  "Type" column already present on the "Locations" sheet)

$wb = $excel.ActiveWorkbook

$wsEvents = $wb.Worksheets.Item("Events")
$wsDestinations = $wb.Worksheets.Item("Event Destinations")

# --- Event Destinations: add column D ("Type" = "Destination") ---
$wsDestinations.Range("D1").Value = "Type"
$wsDestinations.Range("D2").Value = "Destination"
$wsDestinations.Range("D3").Value = "Destination"
$wsDestinations.Range("D4").Value = "Destination"
$wsDestinations.Range("D5").Value = "Destination"
$wsDestinations.Range("D6").Value = "Destination"
$wsDestinations.Range("D7").Value = "Destination"
$wsDestinations.Range("D8").Value = "Destination"
$wsDestinations.Range("D9").Value = "Destination"

# Match the header's bold style used by the rest of row 1 (A1:C1)
$wsDestinations.Range("D1").Font.Bold = $true

# --- Update saved cursor / selection state for both affected sheets ---
$wsEvents.Activate()
$wsEvents.Range("E29").Select()

$wsDestinations.Range("D13").Select()

$wsEvents.Activate()
